$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add row 29
$ws.Cells.Item(29, 1).Value = 45996
$ws.Cells.Item(29, 2).Value = 1

# Add row 30
$ws.Cells.Item(30, 1).Value = 45999
$ws.Cells.Item(30, 2).Value = 2

# Match the date format used by the existing date column (A2:A28)
$ws.Range("A2").Copy()
$ws.Range("A29:A30").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update selection to mirror the author's final selection state
$ws.Range("A29:B30").Select()
